$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 305
$ws.Range("C2").Value = 260
$ws.Range("B3").Value = 215
$ws.Range("C3").Value = 90
$ws.Range("B4").Value = 528
$ws.Range("C4").Value = 355
$ws.Range("B5").Value = 565
$ws.Range("C5").Value = 140

$ws.Range("C3").Select()
